# Update the cryptos list with freshly scraped price/volume values.
# Rows 2-51 hold one coin per row with columns:
#   A = rank index, B = coin name, C = link, D = price, E = 1h volume %
#
# Most D-column "prices" that have exactly one decimal separator look like
# genuine numbers to Excel and would otherwise be silently converted from
# text to a numeric value when written through the COM Value setter. Since
# the source data stores these as plain text (sometimes using '.' as a
# thousands separator, e.g. "27.388.74"), we prefix the ones that parse as
# valid numbers with a leading apostrophe so Excel keeps them as text,
# matching the original cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.388.74"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.639.42"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'211.49"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").Value = "'0.531"
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'22.98"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "'0.0610"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "1.872.28"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.647.00"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'64.37"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "27.370.57"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'229.21"
$ws.Range("E18").Value = "  -5.60%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'4.31"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "'9.56"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'146.93"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "'6.96"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'15.52"
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").Value = "'3.11"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "1.410.37"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'0.564"
$ws.Range("D38").Value = "'0.880"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("E41").Value = "  -0.02%  "

# Rows 42/43 swap: FraxShare (was row 42) and mCoin (was row 43) trade places.
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").Value = "'2.47"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.50"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").Value = "'2.23"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  -7.12%  "
$ws.Range("D47").Value = "1.781.02"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("D48").Value = "'1.65"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").Value = "'87.59"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "'0.0986"
$ws.Range("E51").Value = "  -4.09%  "
